$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (was "1BDS"/102954) to "2BADM"/40
$ws.Range("A2").Value = "2BADM"
$ws.Range("B2").Value = 40

# Add new row 3: "3ADS"/38
$ws.Range("A3").Value = "3ADS"
$ws.Range("B3").Value = 38

# Add new row 4: "1BDS"/36 (moved down from row 2)
$ws.Range("A4").Value = "1BDS"
$ws.Range("B4").Value = 36
